# Resolve the "Test Cases" worksheet and make it active, mirroring the
# manual edit captured in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")
$ws.Activate()

# Row 4's Page Object column now points at the INPUT-specific xpath
# instead of the generic wildcard xpath (new shared string).
$ws.Range("C4").Value = ".//INPUT[@id='landingAmount']"

# The "Validate Mortgage PageN" sample rows (29-38) were cleared out.
# Selecting the block first (so the resulting view/selection state mirrors
# what a user would see after the edit) and then clearing contents removes
# the unstyled cells entirely while leaving the styled-but-now-blank cells
# (B32:D38, styles 7/10) behind - exactly like Excel's own Clear Contents.
$rng = $ws.Range("A29:D38")
$rng.Select()
$rng.ClearContents()

# Park the view roughly where the edited block is, like the saved workbook.
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
